# Natmi following Dr Hou advice
# Adds the "ECs" sending/target cluster to the Dhh-Cdon ligand-receptor table,
# expanding the cluster combinations from 2x3 to 3x3 (rows 2-10) and refreshing
# the corresponding specificity/expression statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Dhh/Cdon)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dhh"
$ws.Range("C2").Value = "Cdon"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.853892
$ws.Range("H2").Value = 5.561676
$ws.Range("I2").Value = 0.5711238486747862
$ws.Range("J2").Value = 0.571123848674786
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.498163333333333
$ws.Range("N2").Value = 4.49449
$ws.Range("O2").Value = 0.02101839619520399
$ws.Range("P2").Value = 0.021018396195204
$ws.Range("Q2").Value = 2.77743301836
$ws.Range("R2").Value = 24.99689716524
$ws.Range("S2").Value = 0.01200410732797639
$ws.Range("T2").Value = 0.01200410732797639

# Row 3: ECs -> FAPs (Dhh/Cdon)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dhh"
$ws.Range("C3").Value = "Cdon"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.853892
$ws.Range("H3").Value = 5.561676
$ws.Range("I3").Value = 0.5711238486747862
$ws.Range("J3").Value = 0.571123848674786
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.96588933333333
$ws.Range("N3").Value = 83.897668
$ws.Range("O3").Value = 0.3923458336491322
$ws.Range("P3").Value = 0.3923458336491322
$ws.Range("Q3").Value = 51.845738507952
$ws.Range("R3").Value = 466.611646571568
$ws.Range("S3").Value = 0.2240780625252098
$ws.Range("T3").Value = 0.2240780625252097

# Row 4: ECs -> sCs (Dhh/Cdon)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dhh"
$ws.Range("C4").Value = "Cdon"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.853892
$ws.Range("H4").Value = 5.561676
$ws.Range("I4").Value = 0.5711238486747862
$ws.Range("J4").Value = 0.571123848674786
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 41.81461766666666
$ws.Range("N4").Value = 125.443853
$ws.Range("O4").Value = 0.5866357701556637
$ws.Range("P4").Value = 0.5866357701556638
$ws.Range("Q4").Value = 77.519785175292
$ws.Range("R4").Value = 697.678066577628
$ws.Range("S4").Value = 0.3350416788215999
$ws.Range("T4").Value = 0.3350416788215999

# Row 5: FAPs -> ECs (Dhh/Cdon)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Dhh"
$ws.Range("C5").Value = "Cdon"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4262446666666667
$ws.Range("H5").Value = 1.278734
$ws.Range("I5").Value = 0.1313121230922664
$ws.Range("J5").Value = 0.1313121230922664
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.498163333333333
$ws.Range("N5").Value = 4.49449
$ws.Range("O5").Value = 0.02101839619520399
$ws.Range("P5").Value = 0.021018396195204
$ws.Range("Q5").Value = 0.6385841306288889
$ws.Range("R5").Value = 5.74725717566
$ws.Range("S5").Value = 0.002759970228386651
$ws.Range("T5").Value = 0.002759970228386651

# Row 6: FAPs -> FAPs (Dhh/Cdon)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dhh"
$ws.Range("C6").Value = "Cdon"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4262446666666667
$ws.Range("H6").Value = 1.278734
$ws.Range("I6").Value = 0.1313121230922664
$ws.Range("J6").Value = 0.1313121230922664
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 27.96588933333333
$ws.Range("N6").Value = 83.897668
$ws.Range("O6").Value = 0.3923458336491322
$ws.Range("P6").Value = 0.3923458336491322
$ws.Range("Q6").Value = 11.92031117692355
$ws.Range("R6").Value = 107.282800592312
$ws.Range("S6").Value = 0.05151976440287273
$ws.Range("T6").Value = 0.05151976440287272

# Row 7: FAPs -> sCs (Dhh/Cdon)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dhh"
$ws.Range("C7").Value = "Cdon"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4262446666666667
$ws.Range("H7").Value = 1.278734
$ws.Range("I7").Value = 0.1313121230922664
$ws.Range("J7").Value = 0.1313121230922664
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 41.81461766666666
$ws.Range("N7").Value = 125.443853
$ws.Range("O7").Value = 0.5866357701556637
$ws.Range("P7").Value = 0.5866357701556638
$ws.Range("Q7").Value = 17.82325776912244
$ws.Range("R7").Value = 160.409319922102
$ws.Range("S7").Value = 0.07703238846100702
$ws.Range("T7").Value = 0.07703238846100702

# Row 8: sCs -> ECs (Dhh/Cdon)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Dhh"
$ws.Range("C8").Value = "Cdon"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9659053333333333
$ws.Range("H8").Value = 2.897716
$ws.Range("I8").Value = 0.2975640282329475
$ws.Range("J8").Value = 0.2975640282329475
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.498163333333333
$ws.Range("N8").Value = 4.49449
$ws.Range("O8").Value = 0.02101839619520399
$ws.Range("P8").Value = 0.021018396195204
$ws.Range("Q8").Value = 1.447083953871111
$ws.Range("R8").Value = 13.02375558484
$ws.Range("S8").Value = 0.006254318638840957
$ws.Range("T8").Value = 0.006254318638840958

# Row 9: sCs -> FAPs (Dhh/Cdon)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Dhh"
$ws.Range("C9").Value = "Cdon"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9659053333333333
$ws.Range("H9").Value = 2.897716
$ws.Range("I9").Value = 0.2975640282329475
$ws.Range("J9").Value = 0.2975640282329475
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.96588933333333
$ws.Range("N9").Value = 83.897668
$ws.Range("O9").Value = 0.3923458336491322
$ws.Range("P9").Value = 0.3923458336491322
$ws.Range("Q9").Value = 27.01240165847644
$ws.Range("R9").Value = 243.111614926288
$ws.Range("S9").Value = 0.1167480067210497
$ws.Range("T9").Value = 0.1167480067210497

# Row 10: sCs -> sCs (Dhh/Cdon)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Dhh"
$ws.Range("C10").Value = "Cdon"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9659053333333333
$ws.Range("H10").Value = 2.897716
$ws.Range("I10").Value = 0.2975640282329475
$ws.Range("J10").Value = 0.2975640282329475
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 41.81461766666666
$ws.Range("N10").Value = 125.443853
$ws.Range("O10").Value = 0.5866357701556637
$ws.Range("P10").Value = 0.5866357701556638
$ws.Range("Q10").Value = 40.38896221552755
$ws.Range("R10").Value = 363.5006599397479
$ws.Range("S10").Value = 0.1745617028730568
$ws.Range("T10").Value = 0.1745617028730568

